$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column C (rows 2-8)
$ws.Range("C2").Value = 1.1002896
$ws.Range("C3").Value = 1.2977649
$ws.Range("C4").Value = 1.4938668
$ws.Range("C5").Value = 1.6867314
$ws.Range("C6").Value = 1.8820485
$ws.Range("C7").Value = 2.0781504
$ws.Range("C8").Value = 2.2744485

# Updated values for column E (rows 2-8) - all set to the same new constant
$ws.Range("E2:E8").Value = 0.01045024556738527
